# Applies the "adhesion" sheet addition + misc tweaks described by the diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "adhesion" worksheet right after the existing "data" sheet.
# ---------------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("data")
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "adhesion"

# ---------------------------------------------------------------------------
# 2. Header row.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Underwater: Bubble"
$ws.Range("B1").Value = "Air"
$ws.Range("C1").Value = "Underwater: Wet"

# ---------------------------------------------------------------------------
# 3. Bulk data block, rows 2-42, columns A, B, C (some cells intentionally blank).
# ---------------------------------------------------------------------------
$A = @(-0.18351022207197126,-0.49231958420420724,-0.67259746581289848,-0.6860405674951997,-0.63794576598504071,-0.58039345613185267,-0.52783122936957105,-0.4830021500754178,-0.44617887710952697,-0.41270886243245025,-0.38495300245023778,-0.36095732692549937,-0.33996676007412885,-0.32145493917399359,-0.30597137177881828,-0.29029541791322144,-0.27717926351627908,-0.26518758982090518,-0.2542330315770604,-0.24383146023970492,-0.23506304300520292,-0.22657219640265719,-0.21921734168810192,-0.21155981355982464,-0.20512396982864983,-0.19864123582428517,-0.19261380726993999,-0.18713426886655876,-0.18187930634311938,-0.17703066111798421,-0.17236876038718646,-0.16788080046034087,-0.16354532472936431,-0.15949225757837102,-0.15560466636971998,-0.15184678021198253,-0.068326159835052538,-0.068356686301740244,-0.06838666880248917,-0.068416072485177948,-0.068444773695119768)
$B = @(-0.11660162781978971,-0.42536457162167757,-0.60559331269693917,-0.6189705643819815,-0.57083378638234084,-0.51324068278280033,-0.46063979444976882,-0.41575859336084869,-0.37889649630073979,-0.34538447835608066,-0.31757777753682598,-0.29353792740494056,-0.27250278613344359,-0.25394610050000171,-0.23841965682258606,-0.222705874850959,-0.20955273519618342,-0.19752077195178705,-0.18653054399209165,-0.17608983077010956,-0.16728364061435705,-0.15875591291703017,-0.15136434446105762,-0.14367070424015069,-0.13719959016356148,-0.13068136744465186,-0.12461904220017998,-0.11910394439240693,-0.11381243929164953,-0.1089315037106077,-0.10423575141160782,-0.099714468436928536,-0.095346398381934913,-0.091260586287944145,-0.087341254655370765,-0.08355176008123491,$null,$null,$null,$null,$null)
$C = @(-0.92254829145588191,-1.4313305078559311,-1.6198822631210408,-1.5319084584213298,-1.3675815088856063,-1.2088137825596583,-1.0758008926812215,-0.96591146794810978,-0.88545386784843927,-0.78929150670781179,-0.73455004594080664,-0.69143913409774871,-0.6294236526596102,-0.58716535417247007,-0.54919744447090546,-0.51514655049444891,-0.48474118013430645,-0.45746439382199272,-0.43266444818357291,-0.4097677888523768,-0.38887175589508088,-0.37091067243721942,-0.35338609035306934,-0.33727907040059851,-0.32192163620692371,-0.30628192750015093,-0.29287143973988355,-0.28542114139895625,-0.26759288562061589,-0.25549623385501274,-0.24604055667509769,-0.25329863748879672,-0.2257354812300601,-0.21669223311704455,-0.20834478946340981,-0.20020780922180911,-0.19237143257555517,-0.18494668250213006,-0.17884488350006206,-0.0096958687905157621,$null)

for ($i = 0; $i -lt $A.Length; $i++) {
    $r = $i + 2
    if ($null -ne $A[$i]) { $ws.Cells.Item($r, 1).Value = $A[$i] }
    if ($null -ne $B[$i]) { $ws.Cells.Item($r, 2).Value = $B[$i] }
    if ($null -ne $C[$i]) { $ws.Cells.Item($r, 3).Value = $C[$i] }
}

# ---------------------------------------------------------------------------
# 4. Summary block, rows 44-49.
# Shared-string table order matters (must match the diff: Adhesion,
# Interfacial tension old, interfacial tension new, New adhesion, Wet/Air,
# " "), so the text values are written in that exact sequence first.
# ---------------------------------------------------------------------------
$ws.Range("A45").Value = "Adhesion"
$ws.Range("A46").Value = "Interfacial tension old"
$ws.Range("A47").Value = "interfacial tension new"
$ws.Range("A48").Value = "New adhesion"
$ws.Range("D44").Value = "Wet/Air"
$ws.Range("D49").Value = " "

$ws.Range("B45").Formula = "=MIN(B2:B42)"
$ws.Range("C45").Formula = "=MIN(C2:C42)"
$ws.Range("D45").Formula = "=C45/B45"

$ws.Range("B46").Value = 24
$ws.Range("C46").Value = 48

$ws.Range("B47").Value = 30
$ws.Range("C47").Value = 40

$ws.Range("B48").Formula = "=B45*B47/B46"
$ws.Range("C48").Formula = "=C45*C47/C46"
$ws.Range("D48").Formula = "=C48/B48"

# ---------------------------------------------------------------------------
# 5. View settings for the new sheet (scrolled down, selection on H43) and
#    make it the active / tabSelected sheet.
# ---------------------------------------------------------------------------
$ws.Range("H43").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 37

$ws.Activate()

# ---------------------------------------------------------------------------
# 6. "data" sheet view tweaks: selection moves to F1, no longer tabSelected.
# ---------------------------------------------------------------------------
$dataSheet.Range("F1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 7. Re-create the (hidden) _xlnm._FilterDatabase defined name that scopes to
#    the "data" sheet, covering A1:AF247.
# ---------------------------------------------------------------------------
$fd = $dataSheet.Names.Add("_xlnm._FilterDatabase", "=data!`$A`$1:`$AF`$247")
$fd.Visible = $false

# Re-activate adhesion sheet so it is the persisted "active" tab.
$ws.Activate()
